$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27, shifting old rows 27-28 down to 28-29
$ws.Rows.Item(27).Insert()

# Fill new row 27 with data (copy of row 28's original "context" fields + new values)
$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C27").Value = "Los Lagos"
$ws.Range("D27").Value = 44509
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = 300000000
$ws.Range("G27").Value = "Espárragos"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 550
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = 1700
$ws.Range("N27").Value = "$/kilo"
$ws.Range("O27").Value = "Provincia de Linares"
$ws.Range("P27").Value = 1700
$ws.Range("Q27").Value = 1
$ws.Range("R27").Value = "Hortaliza"
